$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plasma")

$ws.Range("C2").Value = 113.09999999999999
$ws.Range("C3").Value = 10.5
$ws.Range("C4").Value = 115.3
$ws.Range("C5").Value = 14.390000000000001
$ws.Range("D8").Value = 3.8900000000000001
